# Version 1.1. Fixed bug where i iterated over the wrong variable.
#
# Adds three new upgrade rows (26, 27, 28) to Sheet1:
#   26 -> CYBERNETIC_ENHANCEMENTS, cost 10000000, group 0
#   27 -> LEVEL_1_SENTRY_GUN,      cost 100000,   group 1
#   28 -> TITANIUM_MOUSE_BUTTON,   cost 100000,   group 2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-seed the shared-string table in the same order the original commit
# used (TITANIUM_MOUSE_BUTTON, LEVEL_1_SENTRY_GUN, CYBERNETIC_ENHANCEMENTS)
# so the new <si> entries land at indices 30/31/32, then overwrite the
# scratch cells with the real row 28-30 data below.
$ws.Range("Z1").Value = "TITANIUM_MOUSE_BUTTON"
$ws.Range("Z2").Value = "LEVEL_1_SENTRY_GUN"
$ws.Range("Z3").Value = "CYBERNETIC_ENHANCEMENTS"
$ws.Range("Z1:Z3").Clear()

# Row 28: Id 26, CYBERNETIC_ENHANCEMENTS
$ws.Range("A28").Value = 26
$ws.Range("B28").Value = "CYBERNETIC_ENHANCEMENTS"
$ws.Range("C28").Value = 10000000
$ws.Range("D28").Value = 0

# Row 29: Id 27, LEVEL_1_SENTRY_GUN
$ws.Range("A29").Value = 27
$ws.Range("B29").Value = "LEVEL_1_SENTRY_GUN"
$ws.Range("C29").Value = 100000
$ws.Range("D29").Value = 1

# Row 30: Id 28, TITANIUM_MOUSE_BUTTON
$ws.Range("A30").Value = 28
$ws.Range("B30").Value = "TITANIUM_MOUSE_BUTTON"
$ws.Range("C30").Value = 100000
$ws.Range("D30").Value = 2

# Match the "Delta"/cost column formatting used by the rest of the table
# (right-aligned, wrapped Arial 10 cell style) for the three new cost cells.
$ws.Range("C2").Copy()
$ws.Range("C28:C30").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Reproduce the author's final selection state.
$ws.Range("G28:J33").Select() | Out-Null
